$d = $word.ActiveDocument

$pairs = @(
    @("41×15=", "29×16="),
    @("23×14=", "70×76="),
    @("48×63=", "61×98="),
    @("91×41=", "28×56="),
    @("53×53=", "59×77="),
    @("60×64=", "14×70="),
    @("84×18=", "59×99="),
    @("24×29=", "50×24="),
    @("70×44=", "83×90="),
    @("48×55=", "55×81="),
    @("16×54=", "16×35="),
    @("89×35=", "90×84="),
    @("96×90=", "83×78="),
    @("68×92=", "42×51="),
    @("41×75=", "62×28="),
    @("17×56=", "76×51="),
    @("67×19=", "24×92="),
    @("25×90=", "25×77="),
    @("95×98=", "78×95="),
    @("60×49=", "75×78="),
    @("17×17=", "54×73="),
    @("54×25=", "67×35="),
    @("64×78=", "22×38="),
    @("34×88=", "72×22="),
    @("16×65=", "65×76=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
